$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells are stored as text so values like "1.00" / "20.00" keep
# their trailing zeros instead of being coerced to numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '49.844.75'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.56%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.661.87'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '113.19'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.25%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '327.91'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.40%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.69%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.95%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.93'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -3.50%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.00'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.77%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0820'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.82%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +2.26%  '
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +2.54%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.078.45'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.16%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.655.15'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.40%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.869'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.97%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '49.816.48'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.73'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +3.39%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.93'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.16%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.92%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0₃0954'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.78%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '274.62'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.78%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '69.52'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -4.19%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.47%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.28'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -2.72%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.02%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +1.65%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.10%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -2.58%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '35.10'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -4.91%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '49.53'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.34%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.52'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.44%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.60%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '19.23'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -2.64%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.11%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.97%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.29%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.16'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.90%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '23.52'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +4.44%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '128.04'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +3.22%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0349'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +9.44%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +3.72%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -1.40%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.40%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.066.00'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.95%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.14'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +7.10%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -2.75%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.01'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.26%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.03%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '59.28'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.55%  '
